$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Sample-number text correction: "E7760" -> "E7420" (column G, rows 2-27)
# ---------------------------------------------------------------------
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 7).Value = "E7420"
}

# ---------------------------------------------------------------------
# 2. Re-font column G (rows 2-27) to Arial 11 / general alignment.
#    Build the style once on a scratch cell, then copy/paste the format
#    across the whole range so only a single new font + cellXf is
#    minted instead of one per previously-distinct style.
# ---------------------------------------------------------------------
$scratch = $ws.Range("Z100")
$scratch.Value = "scratch"
$scratch.Font.Size = 11
$scratch.Font.Name = "Arial"
$scratch.Copy()
$ws.Range("G2:G27").PasteSpecial(-4122)
$scratch.Clear()

# ---------------------------------------------------------------------
# 3. Column H (rows 2-27): replace the literal boolean constant with an
#    explicit =FALSE() formula (still evaluates to FALSE).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# ---------------------------------------------------------------------
# 4. Move the active selection from H2:H27 to G2:G27.
# ---------------------------------------------------------------------
[void]$ws.Range("G2:G27").Select()

Write-Output "done"
